$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.1661442006269593
$ws.Range("C2").Value = 0.6018808777429467
$ws.Range("J2").Value = 0.01567398119122257
$ws.Range("P2").Value = 0.1347962382445141
$ws.Range("S2").Value = 0.08150470219435736
$ws.Range("B3").Value = 0.01477832512315271
$ws.Range("C3").Value = 0.04433497536945813
$ws.Range("J3").Value = 0.02955665024630542
$ws.Range("P3").Value = 0.7733990147783252
$ws.Range("S3").Value = 0.1379310344827586
$ws.Range("J4").Value = 0.01851851851851852
$ws.Range("P4").Value = 0.7407407407407407
$ws.Range("S4").Value = 0.2407407407407407
$ws.Range("B6").Value = 0.04954954954954955
$ws.Range("D6").Value = 0.02252252252252252
$ws.Range("E6").Value = 0.004504504504504504
$ws.Range("F6").Value = 0.05855855855855856
$ws.Range("J6").Value = 0.2522522522522522
$ws.Range("O6").Value = 0.009009009009009009
$ws.Range("Q6").Value = 0.1891891891891892
$ws.Range("R6").Value = 0.04954954954954955
$ws.Range("S6").Value = 0.3648648648648649
$ws.Range("B7").Value = 0.1690140845070423
$ws.Range("D7").Value = 0.02112676056338028
$ws.Range("E7").Value = 0.007042253521126761
$ws.Range("F7").Value = 0.09859154929577464
$ws.Range("J7").Value = 0.1056338028169014
$ws.Range("O7").Value = 0.02816901408450704
$ws.Range("Q7").Value = 0.1267605633802817
$ws.Range("R7").Value = 0.04225352112676056
$ws.Range("S7").Value = 0.4014084507042254
$ws.Range("B8").Value = 0.1116279069767442
$ws.Range("D8").Value = 0.02325581395348837
$ws.Range("E8").Value = 0.002325581395348837
$ws.Range("F8").Value = 0.07441860465116279
$ws.Range("J8").Value = 0.09767441860465116
$ws.Range("O8").Value = 0.01162790697674419
$ws.Range("Q8").Value = 0.1790697674418605
$ws.Range("R8").Value = 0.07906976744186046
$ws.Range("S8").Value = 0.4209302325581395
$ws.Range("B9").Value = 0.1120331950207469
$ws.Range("D9").Value = 0.03319502074688797
$ws.Range("F9").Value = 0.04564315352697095
$ws.Range("J9").Value = 0.1037344398340249
$ws.Range("O9").Value = 0.02489626556016597
$ws.Range("Q9").Value = 0.1991701244813278
$ws.Range("R9").Value = 0.07053941908713693
$ws.Range("S9").Value = 0.4107883817427386
$ws.Range("B10").Value = 0.1125377643504532
$ws.Range("D10").Value = 0.02265861027190332
$ws.Range("F10").Value = 0.05891238670694864
$ws.Range("J10").Value = 0.1193353474320242
$ws.Range("O10").Value = 0.01661631419939577
$ws.Range("Q10").Value = 0.2046827794561933
$ws.Range("R10").Value = 0.07401812688821752
$ws.Range("S10").Value = 0.3912386706948641
$ws.Range("G11").Value = 0.131578947368421
$ws.Range("J11").Value = 0.1008771929824561
$ws.Range("K11").Value = 0.1929824561403509
$ws.Range("L11").Value = 0.5614035087719298
$ws.Range("S11").Value = 0.0131578947368421
$ws.Range("G12").Value = 0.6766917293233082
$ws.Range("J12").Value = 0.2180451127819549
$ws.Range("K12").Value = 0.01503759398496241
$ws.Range("L12").Value = 0.03759398496240601
$ws.Range("S12").Value = 0.05263157894736842
$ws.Range("G13").Value = 0.65
$ws.Range("J13").Value = 0.3
$ws.Range("S13").Value = 0.05
$ws.Range("F15").Value = 0.01449275362318841
$ws.Range("H15").Value = 0.1400966183574879
$ws.Range("I15").Value = 0.0966183574879227
$ws.Range("J15").Value = 0.3719806763285024
$ws.Range("K15").Value = 0.04830917874396135
$ws.Range("M15").Value = 0.004830917874396135
$ws.Range("O15").Value = 0.07246376811594203
$ws.Range("S15").Value = 0.251207729468599
$ws.Range("F16").Value = 0.02586206896551724
$ws.Range("H16").Value = 0.1767241379310345
$ws.Range("I16").Value = 0.09051724137931035
$ws.Range("J16").Value = 0.4137931034482759
$ws.Range("K16").Value = 0.09913793103448276
$ws.Range("M16").Value = 0.02586206896551724
$ws.Range("O16").Value = 0.05172413793103448
$ws.Range("S16").Value = 0.1163793103448276
$ws.Range("F17").Value = 0.01552106430155211
$ws.Range("H17").Value = 0.2039911308203991
$ws.Range("I17").Value = 0.1352549889135255
$ws.Range("J17").Value = 0.4301552106430155
$ws.Range("K17").Value = 0.06651884700665188
$ws.Range("M17").Value = 0.01330376940133038
$ws.Range("O17").Value = 0.03547671840354767
$ws.Range("S17").Value = 0.09977827050997783
$ws.Range("F18").Value = 0.01818181818181818
$ws.Range("H18").Value = 0.1878787878787879
$ws.Range("I18").Value = 0.1212121212121212
$ws.Range("J18").Value = 0.4666666666666667
$ws.Range("K18").Value = 0.05454545454545454
$ws.Range("M18").Value = 0.01212121212121212
$ws.Range("O18").Value = 0.05454545454545454
$ws.Range("S18").Value = 0.08484848484848485
$ws.Range("F19").Value = 0.01675977653631285
$ws.Range("H19").Value = 0.1907422186751796
$ws.Range("I19").Value = 0.09577015163607343
$ws.Range("J19").Value = 0.4221867517956903
$ws.Range("K19").Value = 0.08699122106943336
$ws.Range("M19").Value = 0.02154828411811652
$ws.Range("N19").Value = 0.0007980845969672786
$ws.Range("O19").Value = 0.0742218675179569
$ws.Range("S19").Value = 0.09098164405426976
